# Reorganize CA_metrics sheet: drop SASA/max_SASA/Q/theta/conformation/class
# columns, move monosaccharides/motifs up next to binding_score, rename
# sum_SASA -> sasa, recompute flexibility, and add a new boolean
# has_multi_node_motifs column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture values that move/survive before we start overwriting cells.
# (Use Value2 - Value on this host reflects the property accessor itself.)
$monosaccharides2 = $ws.Range("J2").Value2
$motifs2          = $ws.Range("K2").Value2
$sasa2            = $ws.Range("D2").Value2

$monosaccharides3 = $ws.Range("J3").Value2
$motifs3          = $ws.Range("K3").Value2
$sasa3            = $ws.Range("D3").Value2

# Clear the old columns C:L (everything past binding_score) so stale cells
# from the wider layout don't linger past the new G column.
$ws.Range("C1:L3").Clear()

# New headers.
$ws.Range("C1").Value2 = "monosaccharides"
$ws.Range("D1").Value2 = "motifs"
$ws.Range("E1").Value2 = "sasa"
$ws.Range("F1").Value2 = "flexibility"
$ws.Range("G1").Value2 = "has_multi_node_motifs"

# Copy header style (bold, bordered) from A1 onto the new header cells.
$ws.Range("A1").Copy()
$ws.Range("C1:G1").PasteSpecial(-4122)

# Row 2 data.
$ws.Range("C2").Value2 = $monosaccharides2
$ws.Range("D2").Value2 = $motifs2
$ws.Range("E2").Value2 = $sasa2
$ws.Range("F2").Value2 = 21.4689035334009
$ws.Range("G2").Value2 = $true

# Row 3 data.
$ws.Range("C3").Value2 = $monosaccharides3
$ws.Range("D3").Value2 = $motifs3
$ws.Range("E3").Value2 = $sasa3
$ws.Range("F3").Value2 = 19.96344974679717
$ws.Range("G3").Value2 = $true
